# Apply the "notified_human" status update + append the new CUST021-CUST025
# customer rows (repeated 4x) to the customer/database worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) H column: "followed_up" -> "notified_human" for the rows that had it.
# ---------------------------------------------------------------------------
$followedUpRows = @(4, 8, 13, 17, 20)
foreach ($r in $followedUpRows) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Value2 -eq "followed_up") {
        $cell.Value = "notified_human"
    }
}

# ---------------------------------------------------------------------------
# 2) Append 20 new rows (22-41): four repeats of the CUST021-CUST025 block.
# ---------------------------------------------------------------------------
$blockRows = @(
    @("CUST021", "Andrew Martinez", "andrew.m@tech.start", "+91-9876543230", "Machine Learning", 65000, "2025-09-25", "unprocessed", "2025-09-06", "pending", "yes"),
    @("CUST022", "Julia Roberts", "julia.r@creative.co", "+91-9876543231", "Graphic Design", 24000, "2025-09-28", "unprocessed", "2025-09-06", "pending", "yes"),
    @("CUST023", "Chris Evans", "chris.e@enterprise.io", "+91-9876543232", "DevOps", 38000, "2025-10-01", "unprocessed", "2025-09-07", "pending", "yes"),
    @("CUST024", "Emma Watson", "emma.w@digital.net", "+91-9876543233", "Content Strategy", 15000, "2025-10-05", "unprocessed", "2025-09-07", "pending", "yes"),
    @("CUST025", "Robert Lee", "robert.l@cloud.tech", "+91-9876543234", "Blockchain", 72000, "2025-10-08", "unprocessed", "2025-09-08", "pending", "yes")
)

# Columns G (Due Date) and I (Last Contact) hold date-shaped text like
# "2025-09-25" that must stay literal text, not auto-converted to a date
# serial by Excel's input parsing. Flip the cell to text format before
# typing the value, then reset the style back to Normal so no stray
# number-format style is left behind on the cell.
$dateCols = @(7, 9)

$startRow = 22
$blockSize = $blockRows.Length
$repeats = 4

for ($block = 0; $block -lt $repeats; $block++) {
    for ($i = 0; $i -lt $blockSize; $i++) {
        $r = $startRow + ($block * $blockSize) + $i
        $data = $blockRows[$i]
        for ($c = 0; $c -lt $data.Length; $c++) {
            $col = $c + 1
            $cell = $ws.Cells.Item($r, $col)
            if ($dateCols -contains $col) {
                $cell.NumberFormat = "@"
                $cell.Value = $data[$c]
                $cell.Style = "Normal"
            } else {
                $cell.Value = $data[$c]
            }
        }
    }
}
